$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CasesTab) and Row 4 (FilesTab) swap which query text occupies which
# shared-string slot: the "cohort" query (currently referenced by B2) needs to
# end up sharing a slot that sorts before the (rewritten) file query that will
# now live in B4. Capture the current text first.
$casesQuery = $ws.Range("B2").Value()

$filesQueryNew = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis) ,(samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['OSA01'] and demo.breed in ['Unknown'] and samp.summarized_sample_type IN ['Normal Cell Line']
WITH DISTINCT f,  s, c, demo, diag,parent,samp
WITH
        f, c, demo, diag, s,parent,samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,parent,samp,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,parent,samp,
        round(factor * value)/factor AS size
RETURN DISTINCT
       coalesce(f.file_name, '') AS `File Name`,
       coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@

# Use a temporary unique placeholder so the original "CasesTab" string is
# freed up and can be re-inserted in its new position in the shared string
# table, ahead of the updated "FilesTab" query text.
$ws.Range("B2").Value = "__TEMP_PLACEHOLDER__"
$ws.Range("B4").Value = $filesQueryNew
$ws.Range("B2").Value = $casesQuery

# Re-entering the cell values can cause Excel to re-autofit the (wrapped-text)
# row heights; restore the original explicit row heights afterwards.
$ws.Rows(2).RowHeight = 56.25
$ws.Rows(4).RowHeight = 56.25

# Update the selected cell to C4
[void]$ws.Range("C4").Select()
